# Add a "Distribution channel code" column into the stock-obsolescence
# sample sheet, inserted before the existing last column ("Actual item" /
# obsolescence percentage header at column I), pushing that column to J.
#
# Final header row: ... H:Territory | I:Distribution channel code | J:Actual item(header text)
# Data rows: I2="TR", I3="GO"; old column-I values (70/20) now live in J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I; Excel's default Insert() shifts existing
# column I (and everything to its right) one position to the right.
$ws.Columns("I:I").Insert()

# Populate the newly inserted column.
$ws.Range("I1").Value = "Distribution channel code"
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"

# Match the slightly narrower width the source workbook ends up with for
# the new column I (the old column I's width/bestFit now belongs to J,
# which Insert() already preserved automatically).
$ws.Columns("I:I").ColumnWidth = 21.65
